$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 155.75
$ws.Range("I42").Value = 7.6666665
$ws.Range("J42").Value = 600
$ws.Range("K42").Value = 22.9999995
$ws.Range("L42").Value = 1800
$ws.Range("M42").Value = 207.0000005
$ws.Range("N42").Value = -2260
$ws.Range("H95").Value = 29500
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 29500
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 29500
$ws.Range("N95").Value = -34992
$ws.Range("H105").Value = 43835.5
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 43835.5
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 43835.5
$ws.Range("N105").Value = -50823.5
$ws.Range("H106").Value = 2666.6667
$ws.Range("I106").Value = 2000
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 2000
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -1369
$ws.Range("N106").Value = -4262
$ws.Range("H132").Value = 2684.5386
$ws.Range("I132").Value = 2100.0195
$ws.Range("J132").Value = 4813.857
$ws.Range("K132").Value = 6300.058499999999
$ws.Range("L132").Value = 14441.571
$ws.Range("M132").Value = -3770.058499999999
$ws.Range("N132").Value = -19501.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1310.3334
$ws.Range("I61").Value = 1077.2424
$ws.Range("J61").Value = 2592.3333
$ws.Range("K61").Value = 1077.2424
$ws.Range("L61").Value = 2592.3333
$ws.Range("M61").Value = -865.2424000000001
$ws.Range("N61").Value = -3016.3333
$ws.Range("H74").Value = 3858.5625
$ws.Range("I74").Value = 648.5714
$ws.Range("J74").Value = 26328.5
$ws.Range("K74").Value = 648.5714
$ws.Range("L74").Value = 26328.5
$ws.Range("M74").Value = 225.4286
$ws.Range("N74").Value = -28076.5
$ws.Range("H77").Value = 3858.5625
$ws.Range("I77").Value = 648.5714
$ws.Range("J77").Value = 26328.5
$ws.Range("K77").Value = 3242.857
$ws.Range("L77").Value = 131642.5
$ws.Range("M77").Value = 1125.143
$ws.Range("N77").Value = -140378.5
$ws.Range("H132").Value = 8421.179
$ws.Range("I132").Value = 5129.423
$ws.Range("J132").Value = 51214
$ws.Range("K132").Value = 15388.269
$ws.Range("L132").Value = 153642
$ws.Range("M132").Value = -12858.269
$ws.Range("N132").Value = -158702
$ws.Range("H136").Value = 1310.3334
$ws.Range("I136").Value = 1077.2424
$ws.Range("J136").Value = 2592.3333
$ws.Range("K136").Value = 3231.7272
$ws.Range("L136").Value = 7776.999899999999
$ws.Range("M136").Value = -681.7272000000003
$ws.Range("N136").Value = -12876.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 830.1667
$ws.Range("I134").Value = 700.7727
$ws.Range("J134").Value = 2253.5
$ws.Range("K134").Value = 2102.3181
$ws.Range("L134").Value = 6760.5
$ws.Range("M134").Value = 432.6819
$ws.Range("N134").Value = -11830.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 2000
$ws.Range("I50").Value = 2000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 2000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -1375
$ws.Range("H58").Value = 802
$ws.Range("I58").Value = 502.0889
$ws.Range("J58").Value = 1341.84
$ws.Range("K58").Value = 502.0889
$ws.Range("L58").Value = 1341.84
$ws.Range("M58").Value = -299.0889
$ws.Range("N58").Value = -1747.84
$ws.Range("H59").Value = 20000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 20000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290
$ws.Range("H68").Value = 22000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 22000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 22000
$ws.Range("N68").Value = -23498
$ws.Range("H71").Value = 22000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 22000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 66000
$ws.Range("N71").Value = -73488
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H134").Value = 1226.1013
$ws.Range("I134").Value = 1174.5483
$ws.Range("J134").Value = 1414.1177
$ws.Range("K134").Value = 3523.6449
$ws.Range("L134").Value = 4242.3531
$ws.Range("M134").Value = -988.6448999999998
$ws.Range("N134").Value = -9312.3531
$ws.Range("H136").Value = 802
$ws.Range("I136").Value = 502.0889
$ws.Range("J136").Value = 1341.84
$ws.Range("K136").Value = 1506.2667
$ws.Range("L136").Value = 4025.52
$ws.Range("M136").Value = 1043.7333
$ws.Range("N136").Value = -9125.52
$ws.Range("M59").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H93").Value = 8844.888999999999
$ws.Range("I93").Value = 1604
$ws.Range("J93").Value = 9750
$ws.Range("K93").Value = 4812
$ws.Range("L93").Value = 29250
$ws.Range("M93").Value = -2940
$ws.Range("N93").Value = -32994
$ws.Range("H95").Value = 3000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 3000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 9000
$ws.Range("N95").Value = -13118
$ws.Range("H98").Value = 327.13635
$ws.Range("I98").Value = 283.94116
$ws.Range("J98").Value = 474
$ws.Range("K98").Value = 851.82348
$ws.Range("L98").Value = 1422
$ws.Range("M98").Value = 646.17652
$ws.Range("N98").Value = -4418
$ws.Range("H99").Value = 4033
$ws.Range("I99").Value = 4033
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 12099
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -9853
$ws.Range("H100").Value = 3557
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3557
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 10671
$ws.Range("N100").Value = -12293
$ws.Range("H101").Value = 7777.778
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 7777.778
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 23333.334
$ws.Range("N101").Value = -28201.334
$ws.Range("H102").Value = 8000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 8000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 24000
$ws.Range("N102").Value = -28868
$ws.Range("H118").Value = 1333.3334
$ws.Range("I118").Value = 100
$ws.Range("J118").Value = 3800
$ws.Range("K118").Value = 300
$ws.Range("L118").Value = 11400
$ws.Range("M118").Value = 943
$ws.Range("N118").Value = -13886
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3454
$ws.Range("I132").Value = 4010.3076
$ws.Range("J132").Value = 1904.2858
$ws.Range("K132").Value = 12030.9228
$ws.Range("L132").Value = 5712.857400000001
$ws.Range("M132").Value = -9500.9228
$ws.Range("N132").Value = -10772.8574

Write-Output "edits applied"